$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1) Insert a brand-new paragraph right before the "Block Diagram"
#    paragraph. It carries two runs: the plain sentence "Pictures of
#    materials if needed" and a following run holding the horizontal-
#    ellipsis "……", which Word's grammar checker flagged (the run is
#    wrapped in a gramStart/gramEnd proofErr pair), just like the
#    target markup.
# --------------------------------------------------------------------
$blockDiagramPara = $d.Paragraphs.Item(2)
$blockDiagramPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item(2)

$picturesXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Pictures of materials if needed</w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>……</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '</w:p>'
$newPara.Range.InsertXML($picturesXml) | Out-Null

# --------------------------------------------------------------------
# 2) The "_GoBack" bookmark used to sit right after "Flowchart"; it now
#    has to sit right after "Function description (Refer code)" (the
#    end of the document). We reposition it by briefly typing a
#    one-character placeholder at the destination, wrapping a fresh
#    "_GoBack" bookmark around that placeholder (Bookmarks.Add silently
#    relocates the single allowed "_GoBack" bookmark, removing the old
#    one on "Flowchart"), and then deleting the placeholder through the
#    bookmark's own Range so the bookmark collapses exactly in place.
# --------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $lastPara.Range
$endRange.MoveEnd(1, -1) | Out-Null
$endRange.Collapse(0)
$endRange.InsertAfter("Z")
$d.Bookmarks.Add("_GoBack", $endRange) | Out-Null
$goBack = $d.Bookmarks.Item("_GoBack")
$goBackRange = $goBack.Range
$goBackRange.Text = ""

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
